# changelog.docx edit:
#  1. The stray "_GoBack" bookmark that sat in front of the
#     "0.2b - changelog" heading is removed.
#  2. The bullet about the collision-detection optimisation is updated:
#     "dziesięciokrotny" ("tenfold") becomes "dwudziestokrotny"
#     ("twentyfold"), and a new "_GoBack" bookmark is left right after the
#     replaced word (immediately before the trailing period) - this is
#     where Word leaves _GoBack after the last edit a user made.

$d = $word.ActiveDocument

# --- Step 1: drop the old _GoBack bookmark near the heading -------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- Step 2: swap the word inside the optimisation bullet ---------------
$hit = $d.Content
$found = $hit.Find.Execute("dziesięciokrotny", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $hit.Text = "dwudziestokrotny"

    # Force a run break right before "dwudziestokrotny" (matching how Word
    # splits runs at an edit point) by dropping a throw-away bookmark there
    # and immediately deleting it again.
    $wordStart = $d.Range($hit.Start, $hit.Start)
    $d.Bookmarks.Add("zzz_tmp_split", $wordStart)
    $d.Bookmarks("zzz_tmp_split").Delete()

    # Re-create _GoBack collapsed right after "dwudziestokrotny" (i.e.
    # right before the closing period) - this also splits the run there,
    # leaving the period in its own run.
    $wordEnd = $d.Range($hit.End, $hit.End)
    $d.Bookmarks.Add("_GoBack", $wordEnd)
}
